# Applies the 2024-02-25 GitHub Actions cryptos-list refresh to Sheet1.
# Each row below lists only the columns (B=Coin, C=Link, D=Price, E=Volume(1h))
# that actually changed for that row, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values that look like plain numbers (e.g. "0.999", "36.58") would
# otherwise be auto-converted to the Number type by the COM layer; the sheet
# stores every Price/Volume cell as text, so such values are written with a
# leading apostrophe to force Excel to keep (and display) them as text -
# exactly like the other Price cells that already contain literal dots
# (e.g. "51.604.71") and can never be mistaken for numbers.
function Set-TextValue($range, [string]$value) {
    if ($value -match '^-?\d+(\.\d+)?$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

# Row 2
Set-TextValue $ws.Range("D2") '51.604.71'
$ws.Range("E2").Value = '  +1.27%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.020.61'
$ws.Range("E3").Value = '  +2.54%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
Set-TextValue $ws.Range("D5") '378.92'
$ws.Range("E5").Value = '  -0.03%  '

# Row 6
Set-TextValue $ws.Range("D6") '102.99'
$ws.Range("E6").Value = '  +2.01%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.546'
$ws.Range("E7").Value = '  +1.25%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.595'
$ws.Range("E9").Value = '  +2.49%  '

# Row 10
Set-TextValue $ws.Range("D10") '36.58'
$ws.Range("E10").Value = '  +1.30%  '

# Row 11
$ws.Range("E11").Value = '  -0.41%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.0859'
$ws.Range("E12").Value = '  +1.05%  '

# Row 13
Set-TextValue $ws.Range("D13") '3.494.36'
$ws.Range("E13").Value = '  +2.62%  '

# Row 14
Set-TextValue $ws.Range("D14") '18.50'
$ws.Range("E14").Value = '  +1.33%  '

# Row 15
Set-TextValue $ws.Range("D15") '7.74'
$ws.Range("E15").Value = '  +0.62%  '

# Row 16
Set-TextValue $ws.Range("D16") '3.018.28'
$ws.Range("E16").Value = '  +2.47%  '

# Row 17
$ws.Range("E17").Value = '  -1.60%  '

# Row 18
Set-TextValue $ws.Range("D18") '10.50'
$ws.Range("E18").Value = '  -13.29%  '

# Row 19
Set-TextValue $ws.Range("D19") '51.600.84'
$ws.Range("E19").Value = '  +1.35%  '

# Row 20
Set-TextValue $ws.Range("D20") '3.04'
$ws.Range("E20").Value = '  +0.06%  '

# Row 21
Set-TextValue $ws.Range("D21") '12.49'
$ws.Range("E21").Value = '  +0.85%  '

# Row 22
Set-TextValue $ws.Range("D22") '0.0₃0961'
$ws.Range("E22").Value = '  +1.11%  '

# Row 23
Set-TextValue $ws.Range("D23") '69.91'
$ws.Range("E23").Value = '  +0.68%  '

# Row 24
Set-TextValue $ws.Range("D24") '269.01'
$ws.Range("E24").Value = '  +0.95%  '

# Row 25
$ws.Range("E25").Value = '  -3.14%  '

# Row 26
Set-TextValue $ws.Range("D26") '8.29'
$ws.Range("E26").Value = '  +1.60%  '

# Row 27
Set-TextValue $ws.Range("D27") '7.58'
$ws.Range("E27").Value = '  +7.21%  '

# Row 28
$ws.Range("E28").Value = '  +5.46%  '

# Row 29
$ws.Range("E29").Value = '  +0.01%  '

# Row 30
Set-TextValue $ws.Range("D30") '26.25'
$ws.Range("E30").Value = '  +2.54%  '

# Row 31
$ws.Range("E31").Value = '  -0.12%  '

# Row 32
$ws.Range("E32").Value = '  +2.52%  '

# Row 33
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D33") '34.13'
$ws.Range("E33").Value = '  +1.94%  '

# Row 34
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D34") '0.0457'
$ws.Range("E34").Value = '  +5.42%  '

# Row 35
$ws.Range("E35").Value = '  +0.46%  '

# Row 36
$ws.Range("E36").Value = '  +0.27%  '

# Row 37
$ws.Range("E37").Value = '  -0.05%  '

# Row 38
$ws.Range("E38").Value = '  +5.75%  '

# Row 39
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D39") '0.290'
$ws.Range("E39").Value = '  +11.20%  '

# Row 40
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D40") '17.15'
$ws.Range("E40").Value = '  +3.63%  '

# Row 41
$ws.Range("E41").Value = '  +3.76%  '

# Row 42
Set-TextValue $ws.Range("D42") '1.86'
$ws.Range("E42").Value = '  +2.86%  '

# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D43") '127.28'
$ws.Range("E43").Value = '  +5.72%  '

# Row 44
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D44") '0.116'
$ws.Range("E44").Value = '  -0.25%  '

# Row 45
Set-TextValue $ws.Range("D45") '3.73'
$ws.Range("E45").Value = '  +7.94%  '

# Row 46
Set-TextValue $ws.Range("D46") '21.66'
$ws.Range("E46").Value = '  +1.40%  '

# Row 47
Set-TextValue $ws.Range("D47") '2.07'
$ws.Range("E47").Value = '  +3.04%  '

# Row 48
Set-TextValue $ws.Range("D48") '2.40'
$ws.Range("E48").Value = '  +3.39%  '

# Row 49
Set-TextValue $ws.Range("D49") '2.030.14'
$ws.Range("E49").Value = '  +0.92%  '

# Row 50
Set-TextValue $ws.Range("D50") '3.317.74'
$ws.Range("E50").Value = '  +2.40%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.0321'
$ws.Range("E51").Value = '  +1.92%  '
